$wb = $excel.ActiveWorkbook

# --- Typography sheet: set Wildcard Characters for the "Default" typography row ---
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("G4").NumberFormat = "@"
$wsTypo.Range("G4").Value = "0123456789"
$wsTypo.Range("G4").Style = "Normal"

# --- Translation sheet: add new text rows ---
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B5").Value = "SingleUseId2"
$wsTrans.Range("C5").Value = "Default"
$wsTrans.Range("D5").Value = "Left"
$wsTrans.Range("E5").Value = "LTR"
$wsTrans.Range("F5").Value = "Counter: <counter_value>"

$wsTrans.Range("B6").Value = "SingleUseId3"
$wsTrans.Range("C6").Value = "Default"
$wsTrans.Range("D6").Value = "Center"
$wsTrans.Range("E6").Value = "LTR"
$wsTrans.Range("F6").Value = "Count"

$wsTrans.Range("B7").Value = "SingleUseId4"
$wsTrans.Range("C7").Value = "Default"
$wsTrans.Range("D7").Value = "Left"
$wsTrans.Range("E7").Value = "LTR"
$wsTrans.Range("F7").NumberFormat = "@"
$wsTrans.Range("F7").Value = "00"
$wsTrans.Range("F7").Style = "Normal"

# Touch row 8 so it materializes as a (blank) row below the new data, matching
# the author's sheet layout (an empty trailing row), without adding any cell.
$wsTrans.Rows.Item(8).Hidden = $true
$wsTrans.Rows.Item(8).Hidden = $false
